$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-04-10 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-11 Tuesday", 2) | Out-Null

# Update each answer cell in the table by position, preserving run formatting
$t = $d.Tables(1)
$t.Cell(1, 1).Range.Text = "61+19=80"
$t.Cell(1, 2).Range.Text = "37+11=48"
$t.Cell(1, 3).Range.Text = "14+40=54"
$t.Cell(1, 4).Range.Text = "54+22=76"
$t.Cell(1, 5).Range.Text = "94-6=88"
$t.Cell(2, 1).Range.Text = "39+25=64"
$t.Cell(2, 2).Range.Text = "83-18=65"
$t.Cell(2, 3).Range.Text = "58-26=32"
$t.Cell(2, 4).Range.Text = "9+58=67"
$t.Cell(2, 5).Range.Text = "33+27=60"
$t.Cell(3, 1).Range.Text = "24+52=76"
$t.Cell(3, 2).Range.Text = "88-8=80"
$t.Cell(3, 3).Range.Text = "93-18=75"
$t.Cell(3, 4).Range.Text = "14-13=1"
$t.Cell(3, 5).Range.Text = "51-4=47"
$t.Cell(4, 1).Range.Text = "66+25=91"
$t.Cell(4, 2).Range.Text = "5+25=30"
$t.Cell(4, 3).Range.Text = "27+2=29"
$t.Cell(4, 4).Range.Text = "95-48=47"
$t.Cell(4, 5).Range.Text = "0+98=98"
$t.Cell(5, 1).Range.Text = "79+14=93"
$t.Cell(5, 2).Range.Text = "27-17=10"
$t.Cell(5, 3).Range.Text = "57-5=52"
$t.Cell(5, 4).Range.Text = "63+23=86"
$t.Cell(5, 5).Range.Text = "62+15=77"
$t.Cell(6, 1).Range.Text = "35+42=77"
$t.Cell(6, 2).Range.Text = "63+22=85"
$t.Cell(6, 3).Range.Text = "13+73=86"
$t.Cell(6, 4).Range.Text = "83-69=14"
$t.Cell(6, 5).Range.Text = "68-15=53"
$t.Cell(7, 1).Range.Text = "37-2=35"
$t.Cell(7, 2).Range.Text = "45+39=84"
$t.Cell(7, 3).Range.Text = "52-27=25"
$t.Cell(7, 4).Range.Text = "59-58=1"
$t.Cell(7, 5).Range.Text = "82+10=92"
$t.Cell(8, 1).Range.Text = "5+12=17"
$t.Cell(8, 2).Range.Text = "54+37=91"
$t.Cell(8, 3).Range.Text = "23+53=76"
$t.Cell(8, 4).Range.Text = "8+32=40"
$t.Cell(8, 5).Range.Text = "16+9=25"
$t.Cell(9, 1).Range.Text = "74+22=96"
$t.Cell(9, 2).Range.Text = "69+10=79"
$t.Cell(9, 3).Range.Text = "9+20=29"
$t.Cell(9, 4).Range.Text = "9+58=67"
$t.Cell(9, 5).Range.Text = "48-24=24"
$t.Cell(10, 1).Range.Text = "69+27=96"
$t.Cell(10, 2).Range.Text = "91-22=69"
$t.Cell(10, 3).Range.Text = "9+19=28"
$t.Cell(10, 4).Range.Text = "5+10=15"
$t.Cell(10, 5).Range.Text = "7+48=55"
$t.Cell(11, 1).Range.Text = "17+58=75"
$t.Cell(11, 2).Range.Text = "97-56=41"
$t.Cell(11, 3).Range.Text = "51-29=22"
$t.Cell(11, 4).Range.Text = "98-68=30"
$t.Cell(11, 5).Range.Text = "30+38=68"
$t.Cell(12, 1).Range.Text = "68+8=76"
$t.Cell(12, 2).Range.Text = "2+62=64"
$t.Cell(12, 3).Range.Text = "74-19=55"
$t.Cell(12, 4).Range.Text = "13+68=81"
$t.Cell(12, 5).Range.Text = "57-41=16"
$t.Cell(13, 1).Range.Text = "12+74=86"
$t.Cell(13, 2).Range.Text = "64+24=88"
$t.Cell(13, 3).Range.Text = "73-3=70"
$t.Cell(13, 4).Range.Text = "41-36=5"
$t.Cell(13, 5).Range.Text = "73-46=27"
$t.Cell(14, 1).Range.Text = "11-10=1"
$t.Cell(14, 2).Range.Text = "36-15=21"
$t.Cell(14, 3).Range.Text = "19+34=53"
$t.Cell(14, 4).Range.Text = "11+0=11"
$t.Cell(14, 5).Range.Text = "66-2=64"
$t.Cell(15, 1).Range.Text = "68-59=9"
$t.Cell(15, 2).Range.Text = "81-62=19"
$t.Cell(15, 3).Range.Text = "70-24=46"
$t.Cell(15, 4).Range.Text = "48-46=2"
$t.Cell(15, 5).Range.Text = "25+47=72"
$t.Cell(16, 1).Range.Text = "50-11=39"
$t.Cell(16, 2).Range.Text = "26+12=38"
$t.Cell(16, 3).Range.Text = "93-46=47"
$t.Cell(16, 4).Range.Text = "83-41=42"
$t.Cell(16, 5).Range.Text = "67-2=65"
$t.Cell(17, 1).Range.Text = "31-2=29"
$t.Cell(17, 2).Range.Text = "23+1=24"
$t.Cell(17, 3).Range.Text = "37+30=67"
$t.Cell(17, 4).Range.Text = "67+32=99"
$t.Cell(17, 5).Range.Text = "12+24=36"
$t.Cell(18, 1).Range.Text = "77-50=27"
$t.Cell(18, 2).Range.Text = "15+55=70"
$t.Cell(18, 3).Range.Text = "38+1=39"
$t.Cell(18, 4).Range.Text = "9+26=35"
$t.Cell(18, 5).Range.Text = "82-24=58"
$t.Cell(19, 1).Range.Text = "90-20=70"
$t.Cell(19, 2).Range.Text = "71-9=62"
$t.Cell(19, 3).Range.Text = "85-43=42"
$t.Cell(19, 4).Range.Text = "95-87=8"
$t.Cell(19, 5).Range.Text = "59-42=17"
$t.Cell(20, 1).Range.Text = "28+45=73"
$t.Cell(20, 2).Range.Text = "73+3=76"
$t.Cell(20, 3).Range.Text = "98-67=31"
$t.Cell(20, 4).Range.Text = "70-28=42"
$t.Cell(20, 5).Range.Text = "98-73=25"
